# Update benchmark figures on the BENCHMARK sheet (2026-01-30 06:58:55 UTC)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value  = "30,46 TRY - 60,94 TRY - 609,43 TRY"
$ws.Range("G4").Value  = "21,27 TRY - 42,55 TRY - 304,71 TRY"
$ws.Range("G5").Value  = "6,09 TRY - 12,19 TRY - 152,35 TRY"

$ws.Range("G6").Value  = "8.300,01 TL - 76,17 TL"
$ws.Range("H6").Value  = ""

$ws.Range("G8").Value  = "15,23 TRY - 30,47 TRY - 304,71 TRY"
$ws.Range("G9").Value  = "10,63 TRY - 21,27 TRY - 152,35 TRY"
$ws.Range("G10").Value = "3,04 TRY - 6,09 TRY - 76,17 TRY"
$ws.Range("G11").Value = "3,04 TRY - 6,09 TRY - 76,17 TRY"

$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 11.380 TL"
$ws.Range("H13").Value = ""

$ws.Range("G14").Value = "8.300 TL - 6,09 TL"
$ws.Range("H14").Value = ""
